# executive-presentation.pptx edit
#
# The authored diff trims this presales deck down to just 4 slides
# (Title, "Why This Solution?", "Business Value - Financial Impact",
# "Risk Mitigation") and removes a stray bold override from a few table
# cells on the slides that survive. Reproduce that using the normal
# PowerPoint slide-management object model: delete the slides that are
# dropped, move the slides that are kept into their final order, then
# touch up the remaining bold runs.

$p = $ppt.ActivePresentation

# --- 1. Drop the old slides 2-4 (Agenda, Executive Summary, Current
#        State) - their slots get taken over by slides that currently
#        sit later in the deck.
$p.Slides.Item(4).Delete()
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()

# --- 2. Drop "Vision - Future State" and "Solution Overview" which sit
#        right after the title slide now.
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()

# At this point the deck reads:
#   1 Presentation Title
#   2 Why This Solution?
#   3 Business Value - Financial Impact
#   4 Business Value - Strategic Benefits
#   5 Implementation Approach
#   6 Risk Mitigation
#   7 Investment Summary
#   8 Timeline & Milestones
#   9 Success Stories
#  10 Our Partnership Advantage
#  11 Next Steps

# --- 3. Move "Risk Mitigation" up to slide 4 so the deck ends on
#        Title / Why This Solution / Business Value / Risk Mitigation.
$riskIdx = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $titleShape = $p.Slides.Item($i).Shapes.Item(1)
    if ($titleShape.HasTextFrame -and $titleShape.TextFrame.TextRange.Text -like "*Risk Mitigation*") {
        $riskIdx = $i
    }
}
$p.Slides.Item($riskIdx).MoveTo(4)

# --- 4. Everything after slide 4 is no longer part of the deck.
while ($p.Slides.Count -gt 4) {
    $p.Slides.Item($p.Slides.Count).Delete()
}

# --- 5. Strip the leftover bold override from a handful of table
#        cells on the two comparison/metric slides.

# Slide 2 - "Why This Solution?" comparison table, last row
# ("[Current limitation 3]" / "[Our advantage 3]").
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(3).Table
$tbl2.Cell($tbl2.Rows.Count, 1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl2.Cell($tbl2.Rows.Count, 2).Shape.TextFrame.TextRange.Font.Bold = $false

# Slide 3 - "Business Value - Financial Impact" metrics table:
# header row ("Metric" / "Value") and the "ROI" / "[Percentage]" row.
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(3).Table
$tbl3.Cell(1, 1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl3.Cell(1, 2).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl3.Cell($tbl3.Rows.Count, 1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl3.Cell($tbl3.Rows.Count, 2).Shape.TextFrame.TextRange.Font.Bold = $false

# Slide 4 - "Risk Mitigation" table: header row ("Risk" /
# "Mitigation Strategy" / "Success Probability") and the "[Risk 3]" row.
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(3).Table
for ($c = 1; $c -le $tbl4.Columns.Count; $c++) {
    $tbl4.Cell(1, $c).Shape.TextFrame.TextRange.Font.Bold = $false
    $tbl4.Cell($tbl4.Rows.Count, $c).Shape.TextFrame.TextRange.Font.Bold = $false
}
